$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E ("Product Name") to hold the new "Client" column.
$ws.Columns.Item(5).Insert()

# Set header and values for the newly inserted "Client" column.
$ws.Cells.Item(1, 5).Value = "Client"
$ws.Cells.Item(2, 5).Value = "FAMS"
$ws.Cells.Item(3, 5).Value = "FAMS"

# Match the column width of the new column to its left neighbor (Assignee_QA, column D).
$ws.Columns.Item(5).ColumnWidth = $ws.Columns.Item(4).ColumnWidth

# Update the active cell selection as recorded in the edit.
$ws.Range("G18").Select()
